$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 28 (2025Q2) metrics
$ws.Range("C28").Value = 349
$ws.Range("D28").Value = 37
$ws.Range("E28").Value = 312
$ws.Range("F28").Value = 5.763239875389408
